# Apply F-column ("想去人数" / interest-count) updates across sheets,
# matching the commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions) - sheet 1
$wsExpo = $wb.Worksheets.Item(1)
$wsExpo.Range("F5").Value  = 1420
$wsExpo.Range("F7").Value  = 603
$wsExpo.Range("F8").Value  = 123
$wsExpo.Range("F9").Value  = 581
$wsExpo.Range("F10").Value = 29
$wsExpo.Range("F11").Value = 665
$wsExpo.Range("F14").Value = 155
$wsExpo.Range("F15").Value = 226

# 演出 (Performances) - sheet 2
$wsShow = $wb.Worksheets.Item(2)
$wsShow.Range("F11").Value = 7
$wsShow.Range("F12").Value = 193

# 本地生活 (Local life) - sheet 3
$wsLocal = $wb.Worksheets.Item(3)
$wsLocal.Range("F2").Value = 6285
$wsLocal.Range("F3").Value = 777
$wsLocal.Range("F4").Value = 1923

# 全部类型 (All types) - sheet 4, a combined listing of the above
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F2").Value  = 6285
$wsAll.Range("F3").Value  = 777
$wsAll.Range("F4").Value  = 1923
$wsAll.Range("F15").Value = 1420
$wsAll.Range("F19").Value = 7
$wsAll.Range("F20").Value = 603
$wsAll.Range("F21").Value = 193
$wsAll.Range("F22").Value = 123
$wsAll.Range("F23").Value = 582
$wsAll.Range("F24").Value = 29
$wsAll.Range("F26").Value = 665
$wsAll.Range("F31").Value = 155
$wsAll.Range("F37").Value = 226
